$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.518.00"
$ws.Range("E2").Value = "  +2.03%  "

$ws.Range("D3").Value = "1.579.16"
$ws.Range("E3").Value = "  +0.43%  "

$ws.Range("E4").Value = "  +0.64%  "

$cellStyle = $ws.Range("D5").Style
$ws.Range("D5").Value = "'212.42"
$ws.Range("D5").Style = $cellStyle
$ws.Range("E5").Value = "  +0.32%  "

$ws.Range("E6").Value = "  +0.08%  "

$ws.Range("E7").Value = "  +0.62%  "

$cellStyle = $ws.Range("D8").Style
$ws.Range("D8").Value = "'46.89"
$ws.Range("D8").Style = $cellStyle
$ws.Range("E8").Value = "  +8.12%  "

$cellStyle = $ws.Range("D9").Style
$ws.Range("D9").Value = "'23.97"
$ws.Range("D9").Style = $cellStyle
$ws.Range("E9").Value = "  +3.57%  "

$ws.Range("E10").Value = "  -0.48%  "

$ws.Range("E11").Value = "  -0.38%  "

$cellStyle = $ws.Range("D12").Style
$ws.Range("D12").Value = "'0.0880"
$ws.Range("D12").Style = $cellStyle
$ws.Range("E12").Value = "  +0.09%  "

$ws.Range("D13").Value = "1.803.78"
$ws.Range("E13").Value = "  +0.43%  "

$ws.Range("D14").Value = "1.580.52"
$ws.Range("E14").Value = "  +0.52%  "

$ws.Range("E15").Value = "  +0.70%  "

$ws.Range("E16").Value = "  -1.08%  "

$ws.Range("D17").Value = "28.537.91"
$ws.Range("E17").Value = "  +2.20%  "

$cellStyle = $ws.Range("D18").Style
$ws.Range("D18").Value = "'62.44"
$ws.Range("D18").Style = $cellStyle
$ws.Range("E18").Value = "  -1.40%  "

$cellStyle = $ws.Range("D19").Style
$ws.Range("D19").Value = "'229.34"
$ws.Range("D19").Style = $cellStyle
$ws.Range("E19").Value = "  +0.53%  "

$ws.Range("E20").Value = "  -0.04%  "

$ws.Range("E21").Value = "  -1.05%  "

$ws.Range("E22").Value = "  +0.60%  "

$ws.Range("E23").Value = "  -3.68%  "

$cellStyle = $ws.Range("D24").Style
$ws.Range("D24").Value = "'9.17"
$ws.Range("D24").Style = $cellStyle
$ws.Range("E24").Value = "  -1.07%  "

$ws.Range("E25").Value = "  +5.39%  "

$cellStyle = $ws.Range("D26").Style
$ws.Range("D26").Value = "'151.60"
$ws.Range("D26").Style = $cellStyle
$ws.Range("E26").Value = "  +0.07%  "

$cellStyle = $ws.Range("D27").Style
$ws.Range("D27").Value = "'15.04"
$ws.Range("D27").Style = $cellStyle
$ws.Range("E27").Value = "  -0.63%  "

$cellStyle = $ws.Range("D28").Style
$ws.Range("D28").Value = "'6.47"
$ws.Range("D28").Style = $cellStyle
$ws.Range("E28").Value = "  -1.36%  "

$ws.Range("E29").Value = "  -1.04%  "

$ws.Range("E30").Value = "  +0.59%  "

$ws.Range("E31").Value = "  -1.55%  "

$ws.Range("E32").Value = "  -1.29%  "

$ws.Range("E33").Value = "  +0.15%  "

$ws.Range("E34").Value = "  +0.03%  "

$ws.Range("D35").Value = "1.398.29"
$ws.Range("E35").Value = "  -0.83%  "

$ws.Range("E36").Value = "  -1.17%  "

$ws.Range("E37").Value = "  -2.81%  "

$ws.Range("E38").Value = "  +2.46%  "

$ws.Range("E39").Value = "  +6.95%  "

$ws.Range("E40").Value = "  -0.10%  "

$cellStyle = $ws.Range("D41").Style
$ws.Range("D41").Value = "'0.532"
$ws.Range("D41").Style = $cellStyle
$ws.Range("E41").Value = "  -1.01%  "

$ws.Range("E42").Value = "  +0.64%  "

$ws.Range("E43").Value = "  -1.13%  "

$cellStyle = $ws.Range("D44").Style
$ws.Range("D44").Value = "'5.61"
$ws.Range("D44").Style = $cellStyle
$ws.Range("E44").Value = "  +0.32%  "

$ws.Range("E45").Value = "  +2.65%  "

$ws.Range("E46").Value = "  +0.89%  "

$cellStyle = $ws.Range("D47").Style
$ws.Range("D47").Value = "'62.70"
$ws.Range("D47").Style = $cellStyle
$ws.Range("E47").Value = "  -1.14%  "

$ws.Range("D48").Value = "1.716.19"
$ws.Range("E48").Value = "  +0.58%  "

$cellStyle = $ws.Range("D49").Style
$ws.Range("D49").Value = "'85.92"
$ws.Range("D49").Style = $cellStyle
$ws.Range("E49").Value = "  -0.84%  "

$ws.Range("E50").Value = "  -2.27%  "

$ws.Range("E51").Value = "  -1.29%  "
